$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Structural changes: remove "Bomat Courier" row, add six new rows for the
#    new cards in the burn deck (Electrostatic Infantry, Kumano Faces
#    Kakkazan, a second Lightning Strike printing, Phoenix Chick and two
#    Young Pyromancer printings).
# ---------------------------------------------------------------------------

# Remove the "Bomat Courier" line (old row 2).
$ws.Rows.Item(2).Delete()

# Insert the new rows at their final (target) row numbers, top-down, so each
# insertion point is already correctly aligned once the earlier insertions
# have been applied. Insert() copies formatting down from the row above.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(17).Insert()

# ---------------------------------------------------------------------------
# 2) Rewrite the full data block (rows 2-22) with the final values, and the
#    grand-total row (23), so the sheet ends up in the exact target state
#    regardless of how the structural shuffle above landed.
# ---------------------------------------------------------------------------

$data = @(
    @("Bonecrusher Giant: Stomp", "Throne of Eldraine", "Normal", 0.71, 4),
    @("Electrostatic Infantry", "Dominaria United", "Normal", 0.18, 1),
    @("Ghitu Lavarunner", "Dominaria", "Normal", 0.06, 4),
    @("Kumano Faces Kakkazan: Etching of Kumano", "Kamigawa: Neon Dynasty", "Normal", 0.16, 1),
    @("Light Up the Stage", "Ravnica Allegiance", "Normal", 0.23, 4),
    @("Lightning Strike", "Core 2019", "Normal", 0.1, 3),
    @("Lightning Strike", "Dominaria United: Promos", "Normal", 0.29, 1),
    @("Monastery Swiftspear", "Khans of Tarkir", "Normal", 0.59, 4),
    @("Phoenix Chick", "Dominaria United", "Normal", 0.25, 2),
    @("Ramunap Ruins", "Hour of Devastation", "Normal", 0.62, 4),
    @("Soul-Scar Mage", "Amonkhet", "Normal", 1.88, 4),
    @("Viashino Pyromancer", "Core 2019", "Normal", 0.1, 2),
    @("Wild Slash", "Magic the Gathering Products", "Normal", 0.23, 4),
    @("Wizard's Lightning", "Dominaria", "Normal", 0.18, 4),
    @("Young Pyromancer", "Double Masters 2022", "Normal", 0.28, 1),
    @("Young Pyromancer", "Magic 2014", "Normal", 0.24, 1),
    @("Chandra, Torch of Defiance", "Magic the Gathering Products", "Normal", 7.18, 1),
    @("Goblin Chainwhirler", "Dominaria", "Normal", 0.51, 3),
    @("Searing Blood", "Magic the Gathering Products", "Normal", 0.36, 4),
    @("Roiling Vortex", "Zendikar Rising", "Normal", 0.63, 4),
    @("Magma Spray", "Amonkhet", "Normal", 0.05, 3)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Formula = "=D$r*E$r"
    $r = $r + 1
}

# Subtotal markers: row 15 & 16 keep an (empty) numeric G cell, row 17 and 22
# carry the actual SUM() subtotals, row 23 is the grand total.
$ws.Range("G15").Value = ""
$ws.Range("G16").Value = ""
$ws.Range("G17").Formula = "=SUM(F2:F17)"
$ws.Range("G22").Formula = "=SUM(F18:F22)"

$ws.Range("A23:E23").ClearContents()
$ws.Range("G23").ClearContents()
$ws.Range("F23").Formula = "=SUM(F2:F22)"

# ---------------------------------------------------------------------------
# 3) Cosmetic changes: row heights, column width and the saved selection.
# ---------------------------------------------------------------------------

$ws.Range("A1:A23").EntireRow.RowHeight = 12.75

# Excel quantises ColumnWidth to whole pixels of the workbook's default font
# (Arial 10 here, ~6px per character) before it is written back out as the
# raw <col width> - 38.98 itself isn't on that pixel grid, so 38.1 is used
# here because it snaps to the closest reachable width (39.0, vs the
# requested 38.98).
$ws.Columns.Item(1).ColumnWidth = 38.1

$ws.Range("A5").Select()
